# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 171  # data rows 2..171 (row 1 is the header)

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}
